$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new value would otherwise be
# auto-converted to a number by Excel (losing the trailing zeros / exact text).
$textCells = @("D5","D6","D11","D16","D19","D20","D24","D25","D29","D37","D39","D44","D45","D46","D47","D48","D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values from the crypto price refresh.
$ws.Range('D2').Value = '26.277.75'
$ws.Range('E2').Value = '  +0.09%  '
$ws.Range('D3').Value = '1.594.51'
$ws.Range('E3').Value = '  +0.31%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').Value = '212.82'
$ws.Range('E5').Value = '  +0.04%  '
$ws.Range('D6').Value = '0.499'
$ws.Range('E6').Value = '  -0.34%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('E8').Value = '  -0.46%  '
$ws.Range('E9').Value = '  -0.30%  '
$ws.Range('E10').Value = '  -1.79%  '
$ws.Range('D11').Value = '0.0852'
$ws.Range('E11').Value = '  +0.29%  '
$ws.Range('D12').Value = '1.817.58'
$ws.Range('E12').Value = '  +0.27%  '
$ws.Range('D13').Value = '1.602.35'
$ws.Range('E13').Value = '  +0.47%  '
$ws.Range('E14').Value = '  -0.40%  '
$ws.Range('E15').Value = '  -2.43%  '
$ws.Range('D16').Value = '63.94'
$ws.Range('E16').Value = '  -0.69%  '
$ws.Range('D17').Value = '26.267.19'
$ws.Range('D18').Value = '0.0₃0722'
$ws.Range('E18').Value = '  -0.46%  '
$ws.Range('D19').Value = '7.45'
$ws.Range('E19').Value = '  +0.07%  '
$ws.Range('D20').Value = '215.55'
$ws.Range('E20').Value = '  +1.16%  '
$ws.Range('E21').Value = '  -0.03%  '
$ws.Range('E22').Value = '  +0.56%  '
$ws.Range('E23').Value = '  -0.20%  '
$ws.Range('D24').Value = '2.11'
$ws.Range('E24').Value = '  -2.31%  '
$ws.Range('D25').Value = '144.78'
$ws.Range('E25').Value = '  +0.07%  '
$ws.Range('E26').Value = '  -0.03%  '
$ws.Range('E27').Value = '  -1.00%  '
$ws.Range('E28').Value = '  +0.86%  '
$ws.Range('D29').Value = '15.13'
$ws.Range('E29').Value = '  -0.18%  '
$ws.Range('E30').Value = '  -0.49%  '
$ws.Range('E31').Value = '  -0.07%  '
$ws.Range('E32').Value = '  -0.09%  '
$ws.Range('D33').Value = '1.429.64'
$ws.Range('E33').Value = '  +6.77%  '
$ws.Range('E34').Value = '  +0.53%  '
$ws.Range('E35').Value = '  -1.04%  '
$ws.Range('E36').Value = '  -0.49%  '
$ws.Range('D37').Value = '0.560'
$ws.Range('E37').Value = '  -5.38%  '
$ws.Range('E38').Value = '  -0.35%  '
$ws.Range('D39').Value = '0.826'
$ws.Range('E39').Value = '  +1.32%  '
$ws.Range('E40').Value = '  +0.90%  '
$ws.Range('E41').Value = '  +0.03%  '
$ws.Range('E42').Value = '  +1.05%  '
$ws.Range('B43').Value = 'RocketPoolETH'
$ws.Range('C43').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D43').Value = '1.730.18'
$ws.Range('E43').Value = '  +0.34%  '
$ws.Range('D44').Value = '0.757'
$ws.Range('E44').Value = '  -0.67%  '
$ws.Range('B45').Value = 'WEMIXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D45').Value = '0.908'
$ws.Range('E45').Value = '  -11.97%  '
$ws.Range('D46').Value = '60.92'
$ws.Range('E46').Value = '  -1.35%  '
$ws.Range('D47').Value = '86.71'
$ws.Range('E47').Value = '  +0.23%  '
$ws.Range('D48').Value = '1.48'
$ws.Range('E48').Value = '  -1.12%  '
$ws.Range('E49').Value = '  -0.43%  '
$ws.Range('D50').Value = '0.0952'
$ws.Range('E50').Value = '  -2.48%  '
$ws.Range('E51').Value = '  +0.03%  '
